# Apply cryptos list update (Thu Jun 27 04:56:50 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.101.76"
$ws.Range("E2").Value = "  -1.04%  "

# Row 3
$ws.Range("D3").Value = "3.383.06"
$ws.Range("E3").Value = "  -0.19%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.88"
$ws.Range("E5").Value = "  -0.89%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.98"
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "3.381.95"
$ws.Range("E8").Value = "  -0.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.469"
$ws.Range("E9").Value = "  -1.25%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.62"
$ws.Range("E10").Value = "  +1.93%  "

# Row 11
$ws.Range("E11").Value = "  -3.01%  "

# Row 12
$ws.Range("E12").Value = "  -2.47%  "

# Row 13
$ws.Range("D13").Value = "3.958.58"
$ws.Range("E13").Value = "  -0.35%  "

# Row 14
$ws.Range("E14").Value = "  +0.63%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.73"
$ws.Range("E15").Value = "  +1.53%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  -2.91%  "

# Row 17
$ws.Range("D17").Value = "3.383.94"
$ws.Range("E17").Value = "  +0.03%  "

# Row 18
$ws.Range("D18").Value = "61.241.16"
$ws.Range("E18").Value = "  -0.98%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.80"
$ws.Range("E19").Value = "  -2.46%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.75"
$ws.Range("E20").Value = "  -1.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.35"
$ws.Range("E21").Value = "  -1.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.33"
$ws.Range("E22").Value = "  -0.80%  "

# Row 23
$ws.Range("D23").Value = "3.521.54"
$ws.Range("E23").Value = "  -0.26%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.550"
$ws.Range("E24").Value = "  -2.24%  "

# Row 25
$ws.Range("E25").Value = "  +0.13%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000124"
$ws.Range("E26").Value = "  -1.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.04"
$ws.Range("E27").Value = "  -0.27%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("E28").Value = "  +12.27%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.64"
$ws.Range("E29").Value = "  -5.22%  "

# Row 30
$ws.Range("E30").Value = "  -0.42%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.42"
$ws.Range("E31").Value = "  -2.47%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").Value = "  -1.72%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.14"
$ws.Range("E33").Value = "  -1.78%  "

# Row 34
$ws.Range("E34").Value = "  -0.05%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.42"
$ws.Range("E35").Value = "  -0.13%  "

# Row 36
$ws.Range("E36").Value = "  -4.26%  "

# Row 37
$ws.Range("E37").Value = "  -1.36%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.82"
$ws.Range("E38").Value = "  -0.78%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.76"
$ws.Range("E39").Value = "  -0.46%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0760"
$ws.Range("E40").Value = "  -3.42%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.53"
$ws.Range("E41").Value = "  +2.05%  "

# Row 42
$ws.Range("E42").Value = "  +0.05%  "

# Row 43
$ws.Range("E43").Value = "  -1.16%  "

# Row 44
$ws.Range("E44").Value = "  -1.72%  "

# Row 45
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.34"
$ws.Range("E45").Value = "  -1.84%  "

# Row 46
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.19"
$ws.Range("E46").Value = "  -3.94%  "

# Row 47
$ws.Range("D47").Value = "2.560.09"
$ws.Range("E47").Value = "  +8.92%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.78"
$ws.Range("E48").Value = "  -1.21%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.90"
$ws.Range("E49").Value = "  +0.04%  "

# Row 50
$ws.Range("E50").Value = "  +3.97%  "

# Row 51
$ws.Range("E51").Value = "  -1.25%  "
